# The commit swaps the two embedded theme parts: the deck's primary theme
# (ppt/theme/theme1.xml, currently "Integral", used by the Slide Master)
# ends up holding the colours that used to live in ppt/theme/theme2.xml
# ("Office Theme", used by the Notes Master), and vice versa.
#
# The font scheme and format scheme (fills/lines/effects) are already
# byte-for-byte identical between the two theme parts in this deck, so the
# only externally-visible difference to reproduce on the reachable theme
# (theme1.xml / the Slide Master's theme) is its 12 colour-scheme entries.

function ToBgrInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Target values: the "Office Theme" colour scheme (previously theme2.xml),
# now becoming theme1.xml's colour scheme. Order matches the COM
# PpColorSchemeIndex layout: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$targetColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $targetColors.Count; $i++) {
    $tcs.Colors($i + 1).RGB = ToBgrInt($targetColors[$i])
}
